$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ENERGY")
$r = $ws.Range("B4")
$r.Font.Color = 0
